$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append a second copy of the "Purchase 22-23" ledger block (rows 1-16) to
# the bottom of the sheet, landing at rows 22-37 (row 21 left blank, mirroring
# the existing gap before row 18 in the original block).
# ---------------------------------------------------------------------------

# Header row (row 1) -> row 22
$ws.Range("A1:F1").Copy($ws.Range("A22"))

# Data + running-total rows (rows 2-16) -> rows 23-37 (value/style copy)
$ws.Range("A2:F16").Copy($ws.Range("A23"))

# Restore the running-total formulas (Copy() above only carries the cached
# values, so re-point each formula at its row-shifted precedents).
$ws.Range("F27").Formula = "=E23+E24+E25+E26+E27-150000"
$ws.Range("F28").Formula = "=F27+E28"
$ws.Range("F29").Formula = "=F28+E29"
$ws.Range("F30").Formula = "=F29-100000"
$ws.Range("F31").Formula = "=F30+E31"
$ws.Range("F32").Formula = "=F31+E32"
$ws.Range("F33").Formula = "=F32-200000"
$ws.Range("F34").Formula = "=F33+E34"
$ws.Range("F35").Formula = "=F34+E35"
$ws.Range("F36").Formula = "=F35+E36"
$ws.Range("F37").Formula = "=F36-400000"

# Match the row heights used by the equivalent rows higher up the sheet.
$ws.Rows("27:37").RowHeight = 14.4

# Update the view so the newly appended rows are where the user left off.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("A38").Select()
